$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.367.16"
$ws.Range("E2").Value = "  +6.75%  "
$ws.Range("D3").Value = "2.621.58"
$ws.Range("E3").Value = "  +9.09%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "509.78"
$c.ClearFormats()
$ws.Range("E5").Value = "  +5.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "156.45"
$c.ClearFormats()
$ws.Range("E6").Value = "  +3.05%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.26%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.605"
$c.ClearFormats()
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "2.665.70"
$ws.Range("E9").Value = "  +10.10%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.27"
$c.ClearFormats()
$ws.Range("E10").Value = "  +11.54%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.ClearFormats()
$ws.Range("E11").Value = "  +5.34%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.346"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "3.104.90"
$ws.Range("E14").Value = "  +9.92%  "
$ws.Range("D15").Value = "60.570.92"
$ws.Range("E15").Value = "  +6.59%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.77"
$c.ClearFormats()
$ws.Range("E16").Value = "  +5.45%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000139"
$c.ClearFormats()
$ws.Range("E17").Value = "  +4.94%  "
$ws.Range("D18").Value = "2.659.06"
$ws.Range("E18").Value = "  +9.68%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.78"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.63%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "351.51"
$c.ClearFormats()
$ws.Range("E20").Value = "  +8.82%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.47"
$c.ClearFormats()
$ws.Range("E21").Value = "  +5.59%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.ClearFormats()
$ws.Range("E22").Value = "  +4.47%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "59.99"
$c.ClearFormats()
$ws.Range("E24").Value = "  +3.81%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.420"
$c.ClearFormats()
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").Value = "2.778.43"
$ws.Range("E26").Value = "  +10.30%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.ClearFormats()
$ws.Range("E27").Value = "  +4.55%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "0.0₃0860"
$ws.Range("E29").Value = "  +10.44%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.50"
$c.ClearFormats()
$ws.Range("E30").Value = "  +3.84%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.01%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "157.14"
$c.ClearFormats()
$ws.Range("E32").Value = "  +5.77%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.48"
$c.ClearFormats()
$ws.Range("E33").Value = "  +5.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.ClearFormats()
$ws.Range("E34").Value = "  +3.92%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.69"
$c.ClearFormats()
$ws.Range("E35").Value = "  +6.28%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.ClearFormats()
$ws.Range("E36").Value = "  +8.83%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.ClearFormats()
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("E38").Value = "  +10.66%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.858"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "304.76"
$c.ClearFormats()
$ws.Range("E40").Value = "  +15.39%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.75"
$c.ClearFormats()
$ws.Range("E41").Value = "  +7.03%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.832"
$c.ClearFormats()
$ws.Range("E42").Value = "  +28.86%  "
$ws.Range("E43").Value = "  +4.08%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.644"
$c.ClearFormats()
$ws.Range("E44").Value = "  +9.30%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0572"
$c.ClearFormats()
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.992"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "19.84"
$c.ClearFormats()
$ws.Range("E48").Value = "  +14.64%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.01"
$c.ClearFormats()
$ws.Range("E49").Value = "  +10.06%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0236"
$c.ClearFormats()
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").Value = "2.030.79"
$ws.Range("E51").Value = "  +9.32%  "
